$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new price records (for Chirimoya, Macroferia Regional de Talca) need to be
# inserted into the weekly log, right above the existing row 138. Insert two
# blank rows there (this shifts the old rows 138:168 down to 140:170, carrying
# over formatting such as the date style on column D).
$ws.Rows("138:139").Insert()

# Row 138: new "Especial" quality entry dated 2023-10-12 (serial 45211)
$ws.Cells.Item(138, 1).Value = 5
$ws.Cells.Item(138, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(138, 3).Value = "Maule"
$ws.Cells.Item(138, 4).Value = 45211
$ws.Cells.Item(138, 5).Value = 7
$ws.Cells.Item(138, 6).Value = "Fruta"
$ws.Cells.Item(138, 7).Value = 100107
$ws.Cells.Item(138, 8).Value = "Otros"
$ws.Cells.Item(138, 9).Value = 100107002
$ws.Cells.Item(138, 10).Value = "Chirimoya"
$ws.Cells.Item(138, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(138, 12).Value = "Especial"
$ws.Cells.Item(138, 13).Value = 100
$ws.Cells.Item(138, 14).Value = 22000
$ws.Cells.Item(138, 15).Value = 22000
$ws.Cells.Item(138, 16).Value = 22000
$ws.Cells.Item(138, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(138, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(138, 19).Value = 2200
$ws.Cells.Item(138, 20).Value = 10

# Row 139: new "Primera" quality entry, same date
$ws.Cells.Item(139, 1).Value = 5
$ws.Cells.Item(139, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(139, 3).Value = "Maule"
$ws.Cells.Item(139, 4).Value = 45211
$ws.Cells.Item(139, 5).Value = 7
$ws.Cells.Item(139, 6).Value = "Fruta"
$ws.Cells.Item(139, 7).Value = 100107
$ws.Cells.Item(139, 8).Value = "Otros"
$ws.Cells.Item(139, 9).Value = 100107002
$ws.Cells.Item(139, 10).Value = "Chirimoya"
$ws.Cells.Item(139, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(139, 12).Value = "Primera"
$ws.Cells.Item(139, 13).Value = 230
$ws.Cells.Item(139, 14).Value = 20000
$ws.Cells.Item(139, 15).Value = 20000
$ws.Cells.Item(139, 16).Value = 20000
$ws.Cells.Item(139, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(139, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(139, 19).Value = 2000
$ws.Cells.Item(139, 20).Value = 10
